# Aggiornamento dati fino al 23 agosto 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44418, 2, 6, 152.5940996948118),
    @(44419, 1, 7, 178.0264496439471),
    @(44420, 0, 6, 152.5940996948118),
    @(44421, 1, 5, 127.1617497456765),
    @(44422, 0, 5, 127.1617497456765),
    @(44423, 0, 4, 101.7293997965412),
    @(44424, 0, 4, 101.7293997965412),
    @(44425, 1, 3, 76.2970498474059),
    @(44426, 0, 2, 50.8646998982706),
    @(44427, 0, 2, 50.8646998982706),
    @(44428, 2, 3, 76.2970498474059),
    @(44429, 0, 3, 76.2970498474059),
    @(44430, 2, 5, 127.1617497456765),
    @(44431, 0, 5, 127.1617497456765)
)

$lastRow = 343

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $lastRow + 1 + $i

    # Copy formatting from the last existing data row so the new row keeps
    # the same cell styles (e.g. the date column style) as the rest of the table.
    $ws.Range("A$lastRow`:D$lastRow").Copy()
    $ws.Range("A$targetRow`:D$targetRow").PasteSpecial(-4122)

    $data = $newRows[$i]
    $ws.Cells.Item($targetRow, 1).Value = $data[0]
    $ws.Cells.Item($targetRow, 2).Value = $data[1]
    $ws.Cells.Item($targetRow, 3).Value = $data[2]
    $ws.Cells.Item($targetRow, 4).Value = $data[3]
}

$excel.CutCopyMode = $false
